$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 228, shifting existing rows 228.. down by one.
$ws.Rows.Item(228).Insert()

# Populate the newly inserted row 228 with the new weekly data point.
$ws.Cells.Item(228, 1).Value = 11
$ws.Cells.Item(228, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(228, 3).Value = "Bíobío"
$ws.Cells.Item(228, 4).Value = 44726
$ws.Cells.Item(228, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(228, 5).Value = 8
$ws.Cells.Item(228, 6).Value = 100114014
$ws.Cells.Item(228, 7).Value = "Betarraga"
$ws.Cells.Item(228, 8).Value = "Sin especificar"
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 650
$ws.Cells.Item(228, 11).Value = 600
$ws.Cells.Item(228, 12).Value = 650
$ws.Cells.Item(228, 13).Value = 627
$ws.Cells.Item(228, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(228, 15).Value = "Región Metropolitana"
$ws.Cells.Item(228, 16).Value = 125
$ws.Cells.Item(228, 17).Value = 5
$ws.Cells.Item(228, 18).Value = "Hortaliza"
